$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new title row above the existing email list.
$ws.Range("A1").EntireRow.Insert()

# The engine's row-insert does not reflow existing hyperlink anchors, so
# drop every hyperlink on the sheet and recreate them pointing at the
# (now shifted down) email cells.
$ws.Range("A1:A10").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:pluradmiles@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:kilometers520@gmail.com")

# Hyperlinks.Add() re-applies formatting; put the two email cells back on
# the workbook's built-in Hyperlink cell style (matches their original look).
$ws.Range("A2:A3").Style = "Hyperlink"

# New header cell.
$ws.Range("A1").Value = "Emails"
$ws.Range("A1").Font.Name = "Calibri"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 14
$ws.Range("A1").Font.Color = 0
$ws.Rows(1).RowHeight = 18

# Widen column A so the header/email text is not clipped.
$ws.Columns("A:A").ColumnWidth = 24.16666666666667

# Match the saved selection/active cell.
$ws.Range("B6").Select()
